# Splendor_cards.xlsx — restructure sheet order & naming
#
# The "nobles" sheet becomes the first tab (and the active/selected tab),
# the remaining tier1/tier2/tier3 sheets keep their relative order after it,
# and every sheet is renamed to its new 0-based position index ("0".."3").
# The hidden _xlnm._FilterDatabase defined name (which pointed at tier1)
# automatically keeps tracking the tier1 sheet/range once it has been moved
# and renamed, since Excel keeps defined names bound to the sheet object
# rather than to its name/position.

$wb = $excel.ActiveWorkbook

# Move the "nobles" worksheet so it becomes the very first tab.
$noblesSheet = $wb.Worksheets.Item("nobles")
$noblesSheet.Move($wb.Worksheets.Item(1))

# Rename every sheet to match its new position (0-based index as a string).
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = [string]($i - 1)
}

# Make the relocated "nobles" sheet (now "0") the active/selected tab,
# with its selection reset to A1 (matching its default, un-selected state).
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()

# Best-effort: restore the recorded window position (not guaranteed to be
# persisted by every COM host, but harmless to attempt).
try {
    $excel.ActiveWindow.Left = 19455
    $excel.ActiveWindow.Top = 2115
} catch {}
